# Update row-1 header labels across worksheets so Power BI can turn the
# first row into a header automatically. Sheets 1-3 and 5 use "Ano"
# (Year) labels, sheet 4 uses "Intervalo" (Interval) labels, and sheet 6
# only has a B1 header cell.

$wb = $excel.ActiveWorkbook

# Worksheets 1, 2, 3: B1..E1 -> "Ano <valor>"
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range("$col" + "1")
        $cell.Value2 = "Ano " + $cell.Value2
    }
}

# Worksheet 4: B1..E1 -> "Intervalo <valor>"
$ws4 = $wb.Worksheets.Item(4)
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws4.Range("$col" + "1")
    $cell.Value2 = "Intervalo " + $cell.Value2
}

# Worksheet 5: B1..E1 -> "Ano <valor>"
$ws5 = $wb.Worksheets.Item(5)
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws5.Range("$col" + "1")
    $cell.Value2 = "Ano " + $cell.Value2
}

# Worksheet 6: only B1 -> "Ano <valor>"
$ws6 = $wb.Worksheets.Item(6)
$cell6 = $ws6.Range("B1")
$cell6.Value2 = "Ano " + $cell6.Value2
